$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "Yes"
$ws.Range("C10").Value = 0.93111638954869302
$ws.Range("D10").Value = "tanh"
$ws.Range("E10").Value = 0.5
$ws.Range("F10").Value = "adam"
$ws.Range("G10").Value = 64

$ws.Range("C7").Select()
